$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "Extra Number"
$ws.Range("G3").Font.Bold = $true
$ws.Range("G4").Value = 123

$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9

$ws.Range("G5").Select()
